$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DQ_Report")
$ws2 = $wb.Worksheets.Item("DQ_Metrics")

# --- Sheet1 (DQ_Report): delete row 20 ---
$ws1.Rows.Item(20).Delete()

# --- Sheet2 (DQ_Metrics): rename headers and insert new columns ---

# Rename C1 and D1
$ws2.Range("C1").Value = "item_completeness_rate"
$ws2.Range("D1").Value = "value_completeness_rate"

# Rename F1
$ws2.Range("F1").Value = "range_plausibility_rate"

# Rename I1
$ws2.Range("I1").Value = "rdCase_dissimilarity_rate"

# Insert 4 new columns after I (J..M) for rdCase_rel_py_ipat, tracerCase_rel_py_ipat,
# unambigous_rdCase_rel_py_ipat, orphaCase_rel_py_ipat
$ws2.Range("J1:M1").EntireColumn.Insert()

$ws2.Range("J1").Value = "rdCase_rel_py_ipat"
$ws2.Range("K1").Value = "tracerCase_rel_py_ipat"
$ws2.Range("L1").Value = "unambigous_rdCase_rel_py_ipat"
$ws2.Range("M1").Value = "orphaCase_rel_py_ipat"
$ws2.Range("J2").Value = 0.17
$ws2.Range("K2").Value = 0.04
$ws2.Range("L2").Value = 0.07000000000000001
$ws2.Range("M2").Value = 0.15

# Insert a new column for tracerCase_no_py after orphaCoding_no_py (now at N) -> becomes O
$ws2.Range("O1:O1").EntireColumn.Insert()
$ws2.Range("O1").Value = "tracerCase_no_py"
$ws2.Range("O2").Value = 4

# Insert 6 new columns at the end for orphaMissing_no_py, implausible_codeLink_no_py,
# missing_item_no_py, missing_value_no_py, outlier_no_py, duplicateCase_no_py
$ws2.Range("V1:AA1").EntireColumn.Insert()
$ws2.Range("V1").Value = "orphaMissing_no_py"
$ws2.Range("W1").Value = "implausible_codeLink_no_py"
$ws2.Range("X1").Value = "missing_item_no_py"
$ws2.Range("Y1").Value = "missing_value_no_py"
$ws2.Range("Z1").Value = "outlier_no_py"
$ws2.Range("AA1").Value = "duplicateCase_no_py"
$ws2.Range("V2").Value = 2
$ws2.Range("W2").Value = 10
$ws2.Range("X2").Value = 0
$ws2.Range("Y2").Value = 0
$ws2.Range("Z2").Value = 5
$ws2.Range("AA2").Value = 1

# Update existing rate/value cells
$ws2.Range("C2").Value = 100
$ws2.Range("D2").Value = 59.51
$ws2.Range("F2").Value = 93.75
$ws2.Range("I2").Value = 96.3
